# Drop the "How to Play" instructions screen strings.
# Rows 28-36 (Key column values UI_INSTRUCTIONS_TITLE .. UI_INSTRUCTIONS_SURFER_DESCRIPTION)
# are removed entirely, shifting all following rows up by 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28:A36").EntireRow.Delete()

$ws.Range("B65").Select()
